$wb = $excel.ActiveWorkbook

# The "Fund Source Data" sheet's D1 header should read "Distribution"
# instead of "Distribution Amount".
$ws = $wb.Worksheets.Item("Fund Source Data")
$ws.Range("D1").Value = "Distribution"
